# Update Laventan interval's upper bound (D12) and the subsequently
# adjusted Colloncuran interval's lower bound (C13) to keep the two
# adjoining time bins aligned at 13.8 Ma.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Laventan: max_ma (D12) 13 -> 13.8
$ws.Range("D12").Value = 13.8

# Colloncuran: min_ma (C13) 14 -> 13.8
$ws.Range("C13").Value = 13.8

# Leave the active selection on C14, matching the saved workbook state.
$ws.Range("C14").Select()
